# Update column F ("dSF") values on Sheet1 to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 0
    4  = 7
    5  = 1
    6  = -1
    7  = -1
    8  = 2
    9  = -5
    10 = 1
    11 = -4
    12 = 3
    13 = -1
    14 = -2
    15 = 9
    17 = -2
    18 = -6
    19 = 3
    21 = 1
    22 = 6
    23 = 3
    24 = -3
    25 = -3
    26 = 2
    28 = 1
    29 = 3
    30 = 9
    31 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
